$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-03-11 Monday", $false, $false, $false, $false, $false, $true, 1, $false, "2024-03-12 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("173×3=519", $false, $false, $false, $false, $false, $true, 1, $false, "913×6=5478", 2) | Out-Null
$d.Content.Find.Execute("447×4=1788", $false, $false, $false, $false, $false, $true, 1, $false, "309×5=1545", 2) | Out-Null
$d.Content.Find.Execute("691×7=4837", $false, $false, $false, $false, $false, $true, 1, $false, "714×7=4998", 2) | Out-Null
$d.Content.Find.Execute("967×8=7736", $false, $false, $false, $false, $false, $true, 1, $false, "425×9=3825", 2) | Out-Null
$d.Content.Find.Execute("747×5=3735", $false, $false, $false, $false, $false, $true, 1, $false, "780×7=5460", 2) | Out-Null
$d.Content.Find.Execute("485×9=4365", $false, $false, $false, $false, $false, $true, 1, $false, "786×4=3144", 2) | Out-Null
$d.Content.Find.Execute("993×9=8937", $false, $false, $false, $false, $false, $true, 1, $false, "429×3=1287", 2) | Out-Null
$d.Content.Find.Execute("998×7=6986", $false, $false, $false, $false, $false, $true, 1, $false, "902×8=7216", 2) | Out-Null
$d.Content.Find.Execute("577×2=1154", $false, $false, $false, $false, $false, $true, 1, $false, "219×9=1971", 2) | Out-Null
$d.Content.Find.Execute("689×2=1378", $false, $false, $false, $false, $false, $true, 1, $false, "209×4=836", 2) | Out-Null
$d.Content.Find.Execute("949×8=7592", $false, $false, $false, $false, $false, $true, 1, $false, "237×2=474", 2) | Out-Null
$d.Content.Find.Execute("628×3=1884", $false, $false, $false, $false, $false, $true, 1, $false, "834×7=5838", 2) | Out-Null
$d.Content.Find.Execute("987×5=4935", $false, $false, $false, $false, $false, $true, 1, $false, "868×5=4340", 2) | Out-Null
$d.Content.Find.Execute("319×6=1914", $false, $false, $false, $false, $false, $true, 1, $false, "250×8=2000", 2) | Out-Null
$d.Content.Find.Execute("369×4=1476", $false, $false, $false, $false, $false, $true, 1, $false, "608×5=3040", 2) | Out-Null
$d.Content.Find.Execute("863×9=7767", $false, $false, $false, $false, $false, $true, 1, $false, "653×2=1306", 2) | Out-Null
$d.Content.Find.Execute("785×7=5495", $false, $false, $false, $false, $false, $true, 1, $false, "259×6=1554", 2) | Out-Null
$d.Content.Find.Execute("429×2=858", $false, $false, $false, $false, $false, $true, 1, $false, "196×3=588", 2) | Out-Null
$d.Content.Find.Execute("554×9=4986", $false, $false, $false, $false, $false, $true, 1, $false, "244×3=732", 2) | Out-Null
$d.Content.Find.Execute("824×3=2472", $false, $false, $false, $false, $false, $true, 1, $false, "686×9=6174", 2) | Out-Null
$d.Content.Find.Execute("968×4=3872", $false, $false, $false, $false, $false, $true, 1, $false, "303×2=606", 2) | Out-Null
$d.Content.Find.Execute("178×9=1602", $false, $false, $false, $false, $false, $true, 1, $false, "112×2=224", 2) | Out-Null
$d.Content.Find.Execute("210×5=1050", $false, $false, $false, $false, $false, $true, 1, $false, "912×7=6384", 2) | Out-Null
$d.Content.Find.Execute("300×8=2400", $false, $false, $false, $false, $false, $true, 1, $false, "388×6=2328", 2) | Out-Null
$d.Content.Find.Execute("377×6=2262", $false, $false, $false, $false, $false, $true, 1, $false, "652×5=3260", 2) | Out-Null
